# This script reproduces the Jan-5-2023 "GitHub Actions" symbol-list refresh
# for the cryptos worksheet: updated Price / Volume(1h) quotes, plus a handful
# of rows (9-15 and 41-42) whose Coin/Link moved because the underlying feed
# re-sorted/re-labelled entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '256.88'
$ws.Range("E2").Value = '-1.03%'

# --- Row 3 ---
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '27.37'
$ws.Range("E3").Value = '-1.73%'

# --- Row 4 ---
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '4.564'
$ws.Range("E4").Value = '-12.57%'

# --- Row 5 ---
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05891'
$ws.Range("E5").Value = '-0.95%'

# --- Row 6 ---
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '6.627'
$ws.Range("E6").Value = '-1.59%'

# --- Row 7 ---
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8577'
$ws.Range("E7").Value = '-1.75%'

# --- Row 8 ---
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9260'
$ws.Range("E8").Value = '-6.38%'

# --- Row 9 ---
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1407'
$ws.Range("E9").Value = '-1.31%'

# --- Row 10 ---
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03576'
$ws.Range("E10").Value = '-1.55%'

# --- Row 11 ---
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07081'
$ws.Range("E11").Value = '-2.48%'

# --- Row 12 ---
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03232'
$ws.Range("E12").Value = '-0.39%'

# --- Row 13 ---
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09207'
$ws.Range("E13").Value = '-0.36%'

# --- Row 14 ---
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001545'
$ws.Range("E14").Value = '-0.26%'

# --- Row 15 ---
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006067'
$ws.Range("E15").Value = '0.49%'

# --- Row 16 ---
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006050'
$ws.Range("E16").Value = '3.14%'

# --- Row 17 ---
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.56%'

# --- Row 18 ---
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '3.199'
$ws.Range("E18").Value = '-1.49%'

# --- Row 19 ---
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '2.202'
$ws.Range("E19").Value = '-0.36%'

# --- Row 20 ---
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3107'
$ws.Range("E20").Value = '-2.06%'

# --- Row 21 ---
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.62%'

# --- Row 22 ---
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.845'
$ws.Range("E22").Value = '8.61%'

# --- Row 23 ---
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04216'
$ws.Range("E23").Value = '0.97%'

# --- Row 24 ---
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001223'
$ws.Range("E24").Value = '0.58%'

# --- Row 25 ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004297'

# --- Row 26 ---
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001200'
$ws.Range("E26").Value = '0.15%'

# --- Row 27 ---
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-21.94%'

# --- Row 40 ---
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03835'
$ws.Range("E40").Value = '-0.90%'

# --- Row 41 ---
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1103'
$ws.Range("E41").Value = '-0.80%'

# --- Row 42 ---
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003966'
$ws.Range("E42").Value = '-27.10%'

# --- Row 43 ---
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002389'
$ws.Range("E43").Value = '0.49%'

# --- Row 44 ---
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.65%'

# --- Row 45 ---
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005442'
$ws.Range("E45").Value = '0.36%'

# --- Row 46 ---
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.16%'

# --- Row 47 ---
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07196'
$ws.Range("E47").Value = '-15.66%'

# --- Row 48 ---
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1359'
$ws.Range("E48").Value = '6,263.81%'

# --- Row 49 ---
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.16%'

# --- Row 50 ---
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.16%'
